# Revert "Powerpoint writer: consolidate text run nodes."
#
# Splits previously-merged "word + trailing space" text runs back into
# separate runs: one run per word, one run per space. Uses
# TextRange.InsertAfter so each inserted chunk becomes its own <a:r>
# run instead of being folded into the preceding run's text.

function Set-RunChunks {
    param(
        $Shape,
        [string[]]$Chunks
    )

    $tr = $Shape.TextFrame.TextRange
    $tr.Text = $Chunks[0]
    $cur = $tr
    for ($i = 1; $i -lt $Chunks.Count; $i++) {
        $cur = $cur.InsertAfter($Chunks[$i])
    }
}

$p = $ppt.ActivePresentation

Set-RunChunks $p.Slides.Item(1).Shapes.Item("Title 1")  @("Slide", " ", "1", " ", "(Content)")
Set-RunChunks $p.Slides.Item(2).Shapes.Item("Title 1")  @("Slide", " ", "2", " ", "(Content)")
Set-RunChunks $p.Slides.Item(3).Shapes.Item("Title 1")  @("Slide", " ", "3", " ", "(Content)")
Set-RunChunks $p.Slides.Item(4).Shapes.Item("Title 1")  @("Slide", " ", "4", " ", "(Content)")

Set-RunChunks $p.Slides.Item(5).Shapes.Item("Title 1")  @("Slide", " ", "5", " ", "(Two", " ", "Content)")

Set-RunChunks $p.Slides.Item(6).Shapes.Item("Title 1")  @("Slide", " ", "6", " ", "(Two", " ", "Content", " ", "Right)")
Set-RunChunks $p.Slides.Item(6).Shapes.Item("TextBox 3") @("an", " ", "image")

Set-RunChunks $p.Slides.Item(7).Shapes.Item("Title 1")  @("Slide", " ", "7", " ", "(Content", " ", "with", " ", "Caption)")
Set-RunChunks $p.Slides.Item(7).Shapes.Item("TextBox 3") @("An", " ", "image")

Set-RunChunks $p.Slides.Item(8).Shapes.Item("Title 1")  @("Slide", " ", "8", " ", "(Comparison)")
Set-RunChunks $p.Slides.Item(8).Shapes.Item("TextBox 3") @("An", " ", "image")

Set-RunChunks $p.Slides.Item(9).Shapes.Item("Title 1")  @("Slide", " ", "9", " ", "(Content)")
Set-RunChunks $p.Slides.Item(10).Shapes.Item("Title 1") @("Slide", " ", "10", " ", "(Content)")
Set-RunChunks $p.Slides.Item(11).Shapes.Item("Title 1") @("Slide", " ", "11", " ", "(Content)")
Set-RunChunks $p.Slides.Item(12).Shapes.Item("Title 1") @("Slide", " ", "12", " ", "(Content)")
